# Fruta / hortaliza, semanal
# Insert a new weekly price-report row for Membrillo (Mercado Mayorista Lo
# Valledor de Santiago) at row 141, pushing the existing rows 141-146 down
# to 142-147. The worksheet's used range grows from A1:T146 to A1:T147.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 141..146 down one row, leaving a blank row 141 to fill in.
$ws.Rows(141).Insert()

$ws.Range("A141").Value = 6
$ws.Range("B141").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C141").Value = "Metropolitana"
$ws.Range("D141").Value = 45021
$ws.Range("E141").Value = 13
$ws.Range("F141").Value = "Fruta"
$ws.Range("G141").Value = 100104
$ws.Range("H141").Value = "Frutos de pepita"
$ws.Range("I141").Value = 100104003
$ws.Range("J141").Value = "Membrillo"
$ws.Range("K141").Value = "Champion"
$ws.Range("L141").Value = "Primera"
$ws.Range("M141").Value = 20
$ws.Range("N141").Value = 250000
$ws.Range("O141").Value = 250000
$ws.Range("P141").Value = 250000
$ws.Range("Q141").Value = "$/bins (450 kilos)"
$ws.Range("R141").Value = "Región de O'Higgins"
$ws.Range("S141").Value = 556
$ws.Range("T141").Value = 450
